$wb = $excel.ActiveWorkbook

# --- CRpUNL sheet: calibration value updates ---
$wsData = $wb.Worksheets.Item("CRpUNL")
$wsData.Range("B2").Value = 0.04
$wsData.Range("B3").Value = 0.04
$wsData.Range("B4").Value = 0.04
$wsData.Range("B12").Value = 0.001
$wsData.Range("B12").Style = "Normal"

# --- Update selections / active cell on each sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Select()
$wsAbout.Range("A9").Select()

$wsData.Select()
$wsData.Range("E14").Select()

# --- Restore window geometry (best effort) ---
$win = $wb.Windows.Item(1)
$win.Left = 28680
$win.Top = -120
$win.Width = 29040
$win.Height = 17520
